# Apply cyclic rotation of data among rows 2, 3 and 4 of the Artfynd sheet:
#   new row2 <- old row3
#   new row3 <- old row4
#   new row4 <- old row2
# (Columns C, P, Q..Y, AD.. etc. that are identical across the three rows
#  are left untouched; only the columns that actually carry different
#  values per-row are moved.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($r) {
    [PSCustomObject]@{
        A  = $ws.Range("A$r").Value()
        B  = $ws.Range("B$r").Value()
        D  = $ws.Range("D$r").Value()
        E  = $ws.Range("E$r").Value()
        F  = $ws.Range("F$r").Value()
        G  = $ws.Range("G$r").Value()
        H  = $ws.Range("H$r").Value()
        Q  = $ws.Range("Q$r").Value()
        R  = $ws.Range("R$r").Value()
        AC = $ws.Range("AC$r").Value()
    }
}

# Snapshot the "before" values of the three rows first, so that writing
# doesn't clobber data we still need to read.
$row2 = Get-RowData 2
$row3 = Get-RowData 3
$row4 = Get-RowData 4

function Set-RowData($r, $data) {
    $ws.Range("A$r").Value = $data.A
    $ws.Range("B$r").Value = $data.B
    $ws.Range("D$r").Value = $data.D
    $ws.Range("E$r").Value = $data.E
    $ws.Range("F$r").Value = $data.F
    $ws.Range("G$r").Value = $data.G
    $ws.Range("H$r").Value = $data.H
    $ws.Range("Q$r").Value = $data.Q
    $ws.Range("R$r").Value = $data.R
    if ($data.AC -eq $null -or $data.AC -eq "") {
        $ws.Range("AC$r").ClearContents()
    } else {
        $ws.Range("AC$r").Value = $data.AC
    }
}

Set-RowData 2 $row3
Set-RowData 3 $row4
Set-RowData 4 $row2

# The "L" (Kön) column: it was an empty-but-present cell on row 3 and row 4,
# and absent on row 2. After the rotation it should be present (empty) on
# rows 2 and 3, and absent on row 4. Copy an already-blank cell so the new
# blank cell round-trips the same way the original blank cells do.
$ws.Range("K2").Copy($ws.Range("L2"))
# L3 was already present and stays present (row3 <- old row4, which had L4 present) - leave as is.
$ws.Range("L4").ClearContents()
